$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (A2 == 0)
$ws.Range("B2").Value = 0.04752969741821289
$ws.Range("C2").Value = 0.01275396936339973
$ws.Range("D2").Value = 0.02940793037414551
$ws.Range("E2").Value = 0.006698448165590678

# Row 3 (A3 == 1)
$ws.Range("B3").Value = 0.05816531181335449
$ws.Range("C3").Value = 0.004787814655134202
$ws.Range("D3").Value = 0.02771334648132324
$ws.Range("E3").Value = 0.003613983533636935

# Row 4 (A4 == 2)
$ws.Range("B4").Value = 0.2795199871063233
$ws.Range("C4").Value = 0.008223539453831995
$ws.Range("D4").Value = 0.0405156135559082
$ws.Range("E4").Value = 0.002483918768720574

# Row 5 (A5 == 3)
$ws.Range("B5").Value = 0.2359130382537842
$ws.Range("C5").Value = 0.04454885809408572
$ws.Range("D5").Value = 0.0258641242980957
$ws.Range("E5").Value = 0.006311647251569457
